$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rng = $ws.Range("B64:D64")
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(7).Weight = -4138
$rng.Borders.Item(7).Color = 13421772
$rng.Borders.Item(10).LineStyle = 1
$rng.Borders.Item(10).Weight = -4138
$rng.Borders.Item(10).Color = 13421772
Write-Host "done"
